$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the flashcard link cell in row 2 (F2 previously held "/flashcard.html")
$ws.Range("F2").ClearContents()

# Move the "/Light.pdf" value from D7 to E7
$ws.Range("D7").Cut($ws.Range("E7"))

# Update the active selection to match the saved view state
$ws.Range("E13").Select()
